# Update "gh-pages" data workbook to the output generated at 456a3b4.
#
# Sheet 1 ("展览" / Exhibition): only the "想去人数" (F) counters change.
# Sheet 2 ("演出" / Performance): the 2024-07-19 ballet event finished/expired
#   and is removed; the remaining rows shift up one.
# Sheet 3 ("本地生活" / Local life): unchanged (header-only sheet).
# Sheet 4 ("全部类型" / All types, the union of every other sheet): the same
#   2024-07-19 row is removed and the same F-counters are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 - 展览 (Exhibition): bump the "want to go" counters in column F.
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item(1)
$exhibitUpdates = @{
    2 = 5545
    3 = 619
    4 = 12464
    5 = 305
    6 = 618
    7 = 188
    8 = 360
    9 = 1153
    10 = 109
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet 2 - 演出 (Performance): drop the expired 2024-07-19 event (row 2),
# letting the remaining rows shift up, then renumber the index column (A).
# ---------------------------------------------------------------------
$wsPerf = $wb.Worksheets.Item(2)
$wsPerf.Rows.Item(2).Delete()
for ($row = 2; $row -le 5; $row++) {
    $wsPerf.Cells.Item($row, 1).Value = $row - 1
}

# ---------------------------------------------------------------------
# Sheet 4 - 全部类型 (All types): same expired-row removal, then refresh the
# same F counters (now shifted up by one row vs. sheet 1) and renumber A.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Rows.Item(2).Delete()
for ($row = 2; $row -le 14; $row++) {
    $wsAll.Cells.Item($row, 1).Value = $row - 1
}
$allUpdates = @{
    2 = 5545
    3 = 619
    5 = 12464
    6 = 305
    7 = 618
    8 = 188
    11 = 360
    12 = 1153
    14 = 109
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
